$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")
$ws3 = $wb.Worksheets.Item("Hoja3")
$ws4 = $wb.Worksheets.Item("Hoja4")

# --- Hoja1: new review comments in column E ---
$ws1.Range("E6").Value = "x"
$ws1.Range("E8").Value = "No veo la forma de asociar un ticket a un cliente sin que este tenga una estadia, podriamos agregar tickets al cliente?"
$ws1.Range("E10").Value = "La lista de servicios deberia estar en la estadia"
$ws1.Range("E11").Value = "Como resolveriamos el tema de marcar cuando una habitacion esta ocupada o esta libre"

# --- Hoja2: new review comments in column E ---
$ws2.Range("E9").Value = "Mismo problema de asociacion entre cliente y ticket"
$ws2.Range("E10").Value = "Faltaria en el diagrama una entidad que represente"
$ws2.Range("E18").Value = "Se agregaria a la estadia, por eso es que falta la lista de servicios"

# --- Hoja3: new review comments in column E ---
$ws3.Range("E2").Value = "Tabla tareas pendientes???"
$ws3.Range("E3").Value = "Idem"

# --- Hoja4: new review comment in column E, and updated use-case description ---
$ws4.Range("E2").Value = "Tabla ofertas? Como las cargamos(podria ser una pantalla para el admin)?"
$ws4.Range("D13").Value = "Ingresa al sistema un tipo de usuario: encargado, recepcionista,cliente"

# --- Restore the per-sheet cursor/selection from the edit session ---
$ws1.Range("E12").Select() | Out-Null
$ws2.Range("E19").Select() | Out-Null
$ws3.Range("E4").Select() | Out-Null

# Hoja4 ends up being the active (last-focused) sheet/tab
$ws4.Range("A11").Select() | Out-Null
